$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet alternates between "filled" data rows (green fill already
# applied) and "plain" data rows (no fill). Both groups get the same new
# look: Arial font + centered horizontal alignment. We build each look
# once on a scratch cell outside the used range and then paste the
# resulting format onto every row of that group, which keeps the
# generated style table compact (one style per distinct look) instead of
# growing a new style for every single formatting call.

# --- "filled" rows: A2:I2, A4:I4, A6:I6, A8:I8, A10:I10, A12:I12, A14:I14, A16:I16
$tplFilled = $ws.Range("K2")
$tplFilled.Interior.Color = 12379352   # same green fill already used on these rows
$tplFilled.Font.Name = "Arial"
$tplFilled.HorizontalAlignment = -4108 # xlCenter

$filledRows = @(2,4,6,8,10,12,14,16)
foreach ($r in $filledRows) {
    $tplFilled.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)   # xlPasteFormats
}
$tplFilled.Clear()

# --- "plain" rows: A3:I3, A5:I5, A7:I7, A9:I9, A11:I11, A13:I13, A15:I15
$tplPlain = $ws.Range("K1")
$tplPlain.Font.Name = "Arial"
$tplPlain.HorizontalAlignment = -4108  # xlCenter

$plainRows = @(3,5,7,9,11,13,15)
foreach ($r in $plainRows) {
    $tplPlain.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)   # xlPasteFormats
}
$tplPlain.Clear()

$excel.CutCopyMode = 0
